$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 64992; B = "Arthur Ramos";        C = "P&D";                    D = "Viagem de negocios"; E = 1; F = 45106; G = 2970.79 },
    @{ Row = 3;  A = 48408; B = "Clarice Alves";        C = "Financeiro";             D = "Viagem de negocios"; E = 5; F = 45085; G = 6726.71 },
    @{ Row = 4;  A = 39644; B = "Josué Gonçalves";      C = "Financeiro";             D = "Viagem de negocios"; E = 5; F = 45092; G = 3091.11 },
    @{ Row = 5;  A = 40004; B = "Ísis Caldeira";        C = "Financeiro";             D = "Doenca";             E = 5; F = 45092; G = 3788.48 },
    @{ Row = 6;  A = 49740; B = "Felipe Vasconcelos";   C = "Financeiro";             D = "Outros";             E = 6; F = 45091; G = 5833.21 },
    @{ Row = 7;  A = 69620; B = "Olívia Nunes";         C = "Atendimento ao Cliente"; D = "Outros";             E = 8; F = 45099; G = 2950.36 },
    @{ Row = 8;  A = 46771; B = "Zoe Pereira";          C = "P&D";                    D = "Outros";             E = 7; F = 45079; G = 3553.44 },
    @{ Row = 9;  A = 91924; B = "Sarah Vargas";         C = "Operacoes";              D = "Outros";             E = 5; F = 45079; G = 9412.59 },
    @{ Row = 10; A = 62409; B = "Luiz Henrique Lima";   C = "Atendimento ao Cliente"; D = "Problemas pessoais"; E = 3; F = 45090; G = 8267 },
    @{ Row = 11; A = 71279; B = "Gustavo Moraes";       C = "Recursos Humanos";       D = "Viagem de negocios"; E = 7; F = 45089; G = 4035.36 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
